$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "149×6="
$t.Cell(1,2).Range.Text = "362×3="
$t.Cell(1,3).Range.Text = "561×7="
$t.Cell(1,4).Range.Text = "387×2="
$t.Cell(1,5).Range.Text = "678×9="
$t.Cell(5,1).Range.Text = "942×7="
$t.Cell(5,2).Range.Text = "432×6="
$t.Cell(5,3).Range.Text = "934×9="
$t.Cell(5,4).Range.Text = "251×5="
$t.Cell(5,5).Range.Text = "134×6="
$t.Cell(10,1).Range.Text = "695×8="
$t.Cell(10,2).Range.Text = "842×7="
$t.Cell(10,3).Range.Text = "230×2="
$t.Cell(10,4).Range.Text = "264×7="
$t.Cell(10,5).Range.Text = "243×8="
$t.Cell(15,1).Range.Text = "144×2="
$t.Cell(15,2).Range.Text = "625×6="
$t.Cell(15,3).Range.Text = "784×6="
$t.Cell(15,4).Range.Text = "157×5="
$t.Cell(15,5).Range.Text = "408×4="
$t.Cell(20,1).Range.Text = "406×8="
$t.Cell(20,2).Range.Text = "902×5="
$t.Cell(20,3).Range.Text = "517×3="
$t.Cell(20,4).Range.Text = "873×2="
$t.Cell(20,5).Range.Text = "387×6="
